# Portugal Segunda Liga - update of league bases (12-06-2024 23:38)
#
# The underlying source data re-ordered several pairs of adjacent match
# rows (the match id in column B and all of its stats in columns B:AD
# swap places between the two rows), while the running index in column A
# and the shared Div/Date in columns C/D stay put (they are identical for
# both rows in a pair anyway).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of row numbers whose B:AD contents must be swapped.
$rowPairs = @(
    @(88, 89),
    @(186, 187),
    @(230, 231),
    @(243, 244),
    @(296, 297)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B" + $r1 + ":AD" + $r1)
    $rng2 = $ws.Range("B" + $r2 + ":AD" + $r2)

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value = $vals2
    $rng2.Value = $vals1
}
